# Default population of student dropdown
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the file-name noise from the student names in column B
$ws.Range("B2").Value = "GULAM GOS HABIB SHAIKH"
$ws.Range("B3").Value = "HABIBA AFZAL ANSARI"
$ws.Range("B4").Value = "INAAYA MOINUDDIN LUHAR"

# Update the active selection to B4 only
$ws.Range("B4").Select()
